$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D3").Value = -0.0125
$ws.Range("E2:E3").Value = 0.0638
$ws.Range("F2:F3").Value = 0.026
$ws.Range("G2:G3").Value = 0.07811595981555833
$ws.Range("H2:H3").Value = 0.07811595981555833
$ws.Range("I2:I3").Value = 0.07353998281817066
$ws.Range("J2:J3").Value = 0.05634425539612973
$ws.Range("K2:K3").Value = 1179.4
$ws.Range("L2:L3").Value = 0.05169451408734681
$ws.Range("M2:M3").Value = 1004.2
$ws.Range("N2:N3").Value = 0.07111797283324599
$ws.Range("O2:O3").Value = 0.8514498897744616
$ws.Range("P2:P3").Value = 503.7
$ws.Range("Q2:Q3").Value = 0.03567229925921729
$ws.Range("R2:R3").Value = 0.4270815668984229
$ws.Range("S2:S3").Value = 500.5000000000001
$ws.Range("T2:T3").Value = 0.498406691894045
$ws.Range("U2:U3").Value = 8671
$ws.Range("V2:V3").Value = 0.6140847863344712
$ws.Range("W2:W3").Value = 0.07035439670239864
$ws.Range("X2:X3").Value = 0.06940111116996654
$ws.Range("Y2:Y3").Value = 0.0009532855324320966
$ws.Range("Z2:Z3").Value = 1.920195261541051
$ws.Range("AA2:AA3").Value = 0.1081919722267071
$ws.Range("AB2:AB3").Value = 0.0557119928489106
$ws.Range("AC2:AC3").Value = 0.05247997937779649
$ws.Range("AD2:AD3").Value = 5304.8
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 5304.8
$ws.Range("AG2:AG3").Value = -3366.2
$ws.Range("AH2:AH3").Value = 0.2730913770913771
$ws.Range("AI2:AI3").Value = 0.242434213532955
$ws.Range("AJ2:AJ3").Value = -0.313018411753766
$ws.Range("AK2:AK3").Value = -0.2548143886634773
$ws.Range("AL2:AL3").Value = 269.3
$ws.Range("AM2:AM3").Value = 269.3
$ws.Range("AN2:AN3").Value = 3.089754790610985
$ws.Range("AO2:AO3").Value = 6.230226513182324
$ws.Range("AP2:AP3").Value = -1.960626710932494
$ws.Range("AQ2:AQ3").Value = 6.230226513182324
